$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cd80"
$ws.Range("C2").Value = "Cd274"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 30.40041266666667
$ws.Range("H2").Value = 91.201238
$ws.Range("I2").Value = 0.8640973522824783
$ws.Range("J2").Value = 0.8640973522824784
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 17.56563533333333
$ws.Range("N2").Value = 52.696906
$ws.Range("O2").Value = 0.8967254513503751
$ws.Range("P2").Value = 0.8967254513503751
$ws.Range("Q2").Value = 534.0025628855142
$ws.Range("R2").Value = 4806.023065969628
$ws.Range("S2").Value = 0.7748580882361694
$ws.Range("T2").Value = 0.7748580882361695

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cd80"
$ws.Range("C3").Value = "Cd274"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 30.40041266666667
$ws.Range("H3").Value = 91.201238
$ws.Range("I3").Value = 0.8640973522824783
$ws.Range("J3").Value = 0.8640973522824784
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 1.33491
$ws.Range("N3").Value = 4.00473
$ws.Range("O3").Value = 0.06814713783739766
$ws.Range("P3").Value = 0.06814713783739766
$ws.Range("Q3").Value = 40.58181487286
$ws.Range("R3").Value = 365.23633385574
$ws.Range("S3").Value = 0.05888576137092441
$ws.Range("T3").Value = 0.05888576137092442

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Cd80"
$ws.Range("C4").Value = "Cd274"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 30.40041266666667
$ws.Range("H4").Value = 91.201238
$ws.Range("I4").Value = 0.8640973522824783
$ws.Range("J4").Value = 0.8640973522824784
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.6880983333333334
$ws.Range("N4").Value = 2.064295
$ws.Range("O4").Value = 0.03512741081222724
$ws.Range("P4").Value = 0.03512741081222724
$ws.Range("Q4").Value = 20.91847328857889
$ws.Range("R4").Value = 188.26625959721
$ws.Range("S4").Value = 0.03035350267538446
$ws.Range("T4").Value = 0.03035350267538446

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cd80"
$ws.Range("C5").Value = "Cd274"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.356394000000001
$ws.Range("H5").Value = 13.069182
$ws.Range("I5").Value = 0.1238255730991045
$ws.Range("J5").Value = 0.1238255730991045
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 17.56563533333333
$ws.Range("N5").Value = 52.696906
$ws.Range("O5").Value = 0.8967254513503751
$ws.Range("P5").Value = 0.8967254513503751
$ws.Range("Q5").Value = 76.52282837232134
$ws.Range("R5").Value = 688.7054553508921
$ws.Range("S5").Value = 0.1110375429260133
$ws.Range("T5").Value = 0.1110375429260133

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Cd80"
$ws.Range("C6").Value = "Cd274"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.356394000000001
$ws.Range("H6").Value = 13.069182
$ws.Range("I6").Value = 0.1238255730991045
$ws.Range("J6").Value = 0.1238255730991045
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 1.33491
$ws.Range("N6").Value = 4.00473
$ws.Range("O6").Value = 0.06814713783739766
$ws.Range("P6").Value = 0.06814713783739766
$ws.Range("Q6").Value = 5.815393914540001
$ws.Range("R6").Value = 52.33854523086001
$ws.Range("S6").Value = 0.008438358397779434
$ws.Range("T6").Value = 0.008438358397779434

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Cd80"
$ws.Range("C7").Value = "Cd274"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.356394000000001
$ws.Range("H7").Value = 13.069182
$ws.Range("I7").Value = 0.1238255730991045
$ws.Range("J7").Value = 0.1238255730991045
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.6880983333333334
$ws.Range("N7").Value = 2.064295
$ws.Range("O7").Value = 0.03512741081222724
$ws.Range("P7").Value = 0.03512741081222724
$ws.Range("Q7").Value = 2.997627450743334
$ws.Range("R7").Value = 26.97864705669
$ws.Range("S7").Value = 0.004349671775311718
$ws.Range("T7").Value = 0.004349671775311718

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Cd80"
$ws.Range("C8").Value = "Cd274"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.424892
$ws.Range("H8").Value = 1.274676
$ws.Range("I8").Value = 0.01207707461841714
$ws.Range("J8").Value = 0.01207707461841714
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 17.56563533333333
$ws.Range("N8").Value = 52.696906
$ws.Range("O8").Value = 0.8967254513503751
$ws.Range("P8").Value = 0.8967254513503751
$ws.Range("Q8").Value = 7.463497928050668
$ws.Range("R8").Value = 67.171481352456
$ws.Range("S8").Value = 0.01082982018819227
$ws.Range("T8").Value = 0.01082982018819227

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Cd80"
$ws.Range("C9").Value = "Cd274"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.424892
$ws.Range("H9").Value = 1.274676
$ws.Range("I9").Value = 0.01207707461841714
$ws.Range("J9").Value = 0.01207707461841714
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 1.33491
$ws.Range("N9").Value = 4.00473
$ws.Range("O9").Value = 0.06814713783739766
$ws.Range("P9").Value = 0.06814713783739766
$ws.Range("Q9").Value = 0.56719257972
$ws.Range("R9").Value = 5.104733217480001
$ws.Range("S9").Value = 0.0008230180686938094
$ws.Range("T9").Value = 0.0008230180686938094

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Cd80"
$ws.Range("C10").Value = "Cd274"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.424892
$ws.Range("H10").Value = 1.274676
$ws.Range("I10").Value = 0.01207707461841714
$ws.Range("J10").Value = 0.01207707461841714
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.6880983333333334
$ws.Range("N10").Value = 2.064295
$ws.Range("O10").Value = 0.03512741081222724
$ws.Range("P10").Value = 0.03512741081222724
$ws.Range("Q10").Value = 0.2923674770466667
$ws.Range("R10").Value = 2.63130729342
$ws.Range("S10").Value = 0.0004242363615310613
$ws.Range("T10").Value = 0.0004242363615310613
